$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set K5 to the new value "99201 seconds" (new shared string), following the
# pattern of the other cells in column K (Encounter and Measurements Expansion time)
$ws.Range("K5").Value = "99201 seconds"

# Update the sheet's selected/active cell to K8 (matches recorded UI state)
$ws.Range("K8").Select()
